# "finalized translations to element list"
#
# The English translation block (rows 36-66) of the Tabelle1 sheet gets its
# remaining German description texts (column B, rows 49-60) replaced by their
# finished English translations, and the last "element" row (64, formerly
# "Earth" / its German description) is renamed to "Dirt" with a brand new
# English description. Rows 37-48's English descriptions were already final
# and are left untouched; rows 61-63/65/66 keep their English element names
# and only get finished English descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B49").Value = "A weapon is a device which kills people."
$ws.Range("B50").Value = "Fishs are animals which can be fished."
$ws.Range("B51").Value = "Small, flying, bloodsuckers."
$ws.Range("B52").Value = "The death is coming in the future for anybody."
$ws.Range("B53").Value = "You are a human, arent you?"
$ws.Range("B54").Value = "A plant is living in a flowerpot and has dignity and rights."
$ws.Range("B55").Value = "The king is a normal human. But he has a crown on his head."
$ws.Range("B56").Value = "The state similarlly to peace, in which humans are killing them mutually."
$ws.Range("B57").Value = "Multiple humans."
$ws.Range("B58").Value = "Wheat is a plant which grows on the ground… Like literally any other plant."
$ws.Range("B59").Value = "Tree are made of paper and are available in any paper-store."
$ws.Range("B60").Value = "Zombies are people who rose from death."

$ws.Range("B61").Value = "Fire describes the formation of flames during burning."
$ws.Range("B62").Value = "Water is the chemical compound of hydrogen and oxygen."
$ws.Range("B63").Value = "Air is the gas mix of earths atmosphere."

# "Earth" (the dirt/soil element) is renamed to "Dirt" with its own finished
# English description, so it doesn't collide with the "Earth" (the planet)
# sense of the word used elsewhere (e.g. the Sun/Earthquake descriptions).
$ws.Range("A64").Value = "Dirt"
$ws.Range("B64").Value = "Dirt is the death substance which is on the ground."

$ws.Range("B65").Value = "The energy on earth is unsustainable."
$ws.Range("B66").Value = "Intelligence describes the comprehention and thinking of something."

# Reflect where the author ended up working/looking when they finished up.
$win = $excel.ActiveWindow
$ws.Range("B66").Select()
$win.ScrollRow = 39
$win.ScrollColumn = 1
